$d = $word.ActiveDocument

# Locate the paragraph that ends with "Thomson Pioneira (2008)." and delete
# the three paragraphs that used to follow it in the footer block:
#   - an empty paragraph
#   - "Ver no Jupiter Salvar em pdf Salvar em docx"
#   - "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
#      pages. Original theme under Creative Commons Attribution"
# leaving the following empty paragraph (and page-break paragraph) intact.

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Thomson Pioneira (2008).") {
        $target = $i
        break
    }
}

if ($target -ne $null) {
    $startPara = $d.Paragraphs.Item($target + 1)
    $endPara = $d.Paragraphs.Item($target + 3)
    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $r.Delete()
}
